$d = $word.ActiveDocument

# Append a blank spacer paragraph after the current last paragraph
# ("processing/relative_abundance_table.csv"), then append a further
# paragraph containing the status note, matching the other blank-line /
# text pairs already used throughout the document.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$d.Paragraphs.Last.Range.Text = "test ongoing on test branch"
